$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 566.3333
$ws.Range("I2").Value = 499.5
$ws.Range("K2").Value = 499.5
$ws.Range("M2").Value = -386.5
$ws.Range("H9").Value = 233.625
$ws.Range("I9").Value = 264.83334
$ws.Range("K9").Value = 264.83334
$ws.Range("M9").Value = -95.83334000000002
$ws.Range("H33").Value = 185.29411
$ws.Range("I33").Value = 185.29411
$ws.Range("K33").Value = 185.29411
$ws.Range("M33").Value = 43.70589000000001
$ws.Range("H86").Value = 3654
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H87").Value = 24500
$ws.Range("I87").Value = 24500
$ws.Range("K87").Value = 24500
$ws.Range("M87").Value = -23252
$ws.Range("H88").Value = 3597.5557
$ws.Range("I88").Value = 4497
$ws.Range("J88").Value = 3147.8333
$ws.Range("K88").Value = 4497
$ws.Range("L88").Value = 3147.8333
$ws.Range("M88").Value = -4091
$ws.Range("N88").Value = -3959.8333
$ws.Range("H89").Value = 3654
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H90").Value = 24500
$ws.Range("I90").Value = 24500
$ws.Range("K90").Value = 73500
$ws.Range("M90").Value = -67260
$ws.Range("H91").Value = 3597.5557
$ws.Range("I91").Value = 4497
$ws.Range("J91").Value = 3147.8333
$ws.Range("K91").Value = 4497
$ws.Range("L91").Value = 3147.8333
$ws.Range("M91").Value = -3093
$ws.Range("N91").Value = -5955.8333
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -5082
$ws.Range("H137").Value = 2273
$ws.Range("I137").Value = 2292.4
$ws.Range("J137").Value = 2224.5
$ws.Range("K137").Value = 6877.200000000001
$ws.Range("L137").Value = 6673.5
$ws.Range("M137").Value = -4327.200000000001
$ws.Range("N137").Value = -11773.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4390.5454
$ws.Range("I2").Value = 4521
$ws.Range("J2").Value = 4162.25
$ws.Range("K2").Value = 4521
$ws.Range("L2").Value = 4162.25
$ws.Range("M2").Value = -4408
$ws.Range("N2").Value = -4388.25
$ws.Range("H43").Value = 47082.5
$ws.Range("I43").Value = 46342
$ws.Range("K43").Value = 46342
$ws.Range("M43").Value = -46029
$ws.Range("H45").Value = 2734.1875
$ws.Range("I45").Value = 2268
$ws.Range("J45").Value = 5997.5
$ws.Range("K45").Value = 2268
$ws.Range("L45").Value = 5997.5
$ws.Range("M45").Value = -1891
$ws.Range("N45").Value = -6751.5
$ws.Range("H61").Value = 1906.4286
$ws.Range("I61").Value = 1936.75
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 1936.75
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -1724.75
$ws.Range("N61").Value = -1724
$ws.Range("H74").Value = 3646.25
$ws.Range("I74").Value = 3205.3333
$ws.Range("K74").Value = 3205.3333
$ws.Range("M74").Value = -2331.3333
$ws.Range("H77").Value = 3646.25
$ws.Range("I77").Value = 3205.3333
$ws.Range("K77").Value = 16026.6665
$ws.Range("M77").Value = -11658.6665
$ws.Range("H116").Value = 4390.5454
$ws.Range("I116").Value = 4521
$ws.Range("J116").Value = 4162.25
$ws.Range("K116").Value = 4521
$ws.Range("L116").Value = 4162.25
$ws.Range("M116").Value = -2227
$ws.Range("N116").Value = -8750.25
$ws.Range("H136").Value = 1906.4286
$ws.Range("I136").Value = 1936.75
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 5810.25
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -3260.25
$ws.Range("N136").Value = -9000
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4390.5454
$ws.Range("I3").Value = 4521
$ws.Range("J3").Value = 4162.25
$ws.Range("K3").Value = 4521
$ws.Range("L3").Value = 4162.25
$ws.Range("M3").Value = -4407
$ws.Range("N3").Value = -4390.25
$ws.Range("H86").Value = 1866.1666
$ws.Range("I86").Value = 1999.3334
$ws.Range("J86").Value = 1466.6666
$ws.Range("K86").Value = 1999.3334
$ws.Range("L86").Value = 1466.6666
$ws.Range("M86").Value = -876.3334
$ws.Range("N86").Value = -3712.6666
$ws.Range("H89").Value = 1866.1666
$ws.Range("I89").Value = 1999.3334
$ws.Range("J89").Value = 1466.6666
$ws.Range("K89").Value = 9996.666999999999
$ws.Range("L89").Value = 7333.333000000001
$ws.Range("M89").Value = -4380.666999999999
$ws.Range("N89").Value = -18565.333
$ws.Range("H99").Value = 4870.5
$ws.Range("I99").Value = 4976.3076
$ws.Range("K99").Value = 4976.3076
$ws.Range("M99").Value = -3478.3076
$ws.Range("H105").Value = 4744.6
$ws.Range("I105").Value = 4555.75
$ws.Range("K105").Value = 4555.75
$ws.Range("M105").Value = -2808.75
$ws.Range("H134").Value = 6049.1665
$ws.Range("I134").Value = 2195.625
$ws.Range("J134").Value = 13756.25
$ws.Range("K134").Value = 6586.875
$ws.Range("L134").Value = 41268.75
$ws.Range("M134").Value = -4051.875
$ws.Range("N134").Value = -46338.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6106
$ws.Range("I31").Value = 6373.6665
$ws.Range("J31").Value = 4500
$ws.Range("K31").Value = 6373.6665
$ws.Range("L31").Value = 4500
$ws.Range("M31").Value = -6078.6665
$ws.Range("N31").Value = -5090
$ws.Range("H34").Value = 6106
$ws.Range("I34").Value = 6373.6665
$ws.Range("J34").Value = 4500
$ws.Range("K34").Value = 6373.6665
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -6171.6665
$ws.Range("N34").Value = -4904
$ws.Range("H58").Value = 6099
$ws.Range("I58").Value = 6332.222
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 6332.222
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -6129.222
$ws.Range("N58").Value = -4406
$ws.Range("H86").Value = 11399.4
$ws.Range("I86").Value = 8998.5
$ws.Range("J86").Value = 13000
$ws.Range("K86").Value = 8998.5
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = -7875.5
$ws.Range("N86").Value = -15246
$ws.Range("H89").Value = 11399.4
$ws.Range("I89").Value = 8998.5
$ws.Range("J89").Value = 13000
$ws.Range("K89").Value = 44992.5
$ws.Range("L89").Value = 65000
$ws.Range("M89").Value = -39376.5
$ws.Range("N89").Value = -76232
$ws.Range("H132").Value = 2069.182
$ws.Range("I132").Value = 1617.8889
$ws.Range("K132").Value = 4853.6667
$ws.Range("M132").Value = -2323.6667
$ws.Range("H134").Value = 1897
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H136").Value = 6099
$ws.Range("I136").Value = 6332.222
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 18996.666
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -16446.666
$ws.Range("N136").Value = -17100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4111
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null
$ws.Range("H55").Value = 2042.3077
$ws.Range("I55").Value = 1350
$ws.Range("J55").Value = 2408.8235
$ws.Range("K55").Value = 4050
$ws.Range("L55").Value = 7226.470499999999
$ws.Range("M55").Value = -3873
$ws.Range("N55").Value = -7580.470499999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H122").Value = 2972
$ws.Range("I122").Value = 2972
$ws.Range("K122").Value = 8916
$ws.Range("M122").Value = -6466
$ws.Range("H132").Value = 6155
$ws.Range("I132").Value = 6259
$ws.Range("K132").Value = 18777
$ws.Range("M132").Value = -16247
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3654.5334
$ws.Range("I40").Value = 3264
$ws.Range("J40").Value = 5216.6665
$ws.Range("K40").Value = 3264
$ws.Range("L40").Value = 5216.6665
$ws.Range("M40").Value = -3128
$ws.Range("N40").Value = -5488.6665
$ws.Range("H100").Value = 3720.8572
$ws.Range("I100").Value = 4166.1665
$ws.Range("J100").Value = 1049
$ws.Range("K100").Value = 4166.1665
$ws.Range("L100").Value = 1049
$ws.Range("M100").Value = -3625.1665
$ws.Range("N100").Value = -2131
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = $null
$ws.Range("H122").Value = 2103.238
$ws.Range("I122").Value = 1731.6111
$ws.Range("K122").Value = 5194.8333
$ws.Range("M122").Value = -2744.8333
$ws.Range("H126").Value = 1210
$ws.Range("I126").Value = 1210
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3630
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1160
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 2548.6667
$ws.Range("I132").Value = 1188.55
$ws.Range("K132").Value = 3565.65
$ws.Range("M132").Value = -1035.65
$ws.Range("H136").Value = 2824
$ws.Range("I136").Value = 2535.4546
$ws.Range("J136").Value = 5998
$ws.Range("K136").Value = 7606.3638
$ws.Range("L136").Value = 17994
$ws.Range("M136").Value = -5056.3638
$ws.Range("N136").Value = -23094
